$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.91
$ws.Range("I2").Value = 3.9
$ws.Range("U2").Value = 9.5
$ws.Range("X2").Value = 17
$ws.Range("AD2").Value = 201
$ws.Range("AG2").Value = 13

# Row 4
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 2.55
$ws.Range("L4").Value = 1.35
$ws.Range("M4").Value = 2.7
$ws.Range("N4").Value = 2.02
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 1.44
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.82
$ws.Range("S4").Value = 1.78
$ws.Range("T4").Value = 7.6
$ws.Range("U4").Value = 12
$ws.Range("V4").Value = 10
$ws.Range("W4").Value = 27
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 35
$ws.Range("Z4").Value = 8.75
$ws.Range("AA4").Value = 6.3
$ws.Range("AB4").Value = 15.5
$ws.Range("AC4").Value = 80
$ws.Range("AD4").Value = 700
$ws.Range("AE4").Value = 7.6
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 27
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 35

# Row 5
$ws.Range("G5").Value = 7.9
$ws.Range("H5").Value = 5.1
$ws.Range("I5").Value = 1.27
$ws.Range("T5").Value = 20
$ws.Range("U5").Value = 45
$ws.Range("V5").Value = 21
$ws.Range("X5").Value = 65
$ws.Range("Z5").Value = 16
$ws.Range("AA5").Value = 9.25
$ws.Range("AB5").Value = 17
$ws.Range("AE5").Value = 7.2
$ws.Range("AF5").Value = 5.9
$ws.Range("AG5").Value = 7.5
$ws.Range("AH5").Value = 6.9
$ws.Range("AJ5").Value = 20

# Row 7
$ws.Range("G7").Value = 1.47
$ws.Range("H7").Value = 4.2
$ws.Range("I7").Value = 5.2
$ws.Range("T7").Value = 8.25
$ws.Range("U7").Value = 7.5
$ws.Range("V7").Value = 7.1
$ws.Range("W7").Value = 9.5
$ws.Range("X7").Value = 9
$ws.Range("Y7").Value = 15.5
$ws.Range("AA7").Value = 7.8
$ws.Range("AB7").Value = 12
$ws.Range("AC7").Value = 37
$ws.Range("AD7").Value = 200
$ws.Range("AE7").Value = 16.5
$ws.Range("AF7").Value = 30
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 75
$ws.Range("AI7").Value = 35
$ws.Range("AJ7").Value = 30

# Row 8
$ws.Range("J8").Value = 1.06
$ws.Range("K8").Value = 10
$ws.Range("N8").Value = 1.95
$ws.Range("O8").Value = 1.85
$ws.Range("V8").Value = 8.5
$ws.Range("Z8").Value = 10
$ws.Range("AD8").Value = 251
$ws.Range("AE8").Value = 12
$ws.Range("AF8").Value = 19
$ws.Range("AJ8").Value = 41

# Row 14
$ws.Range("U14").Value = 11
$ws.Range("V14").Value = 9
$ws.Range("AF14").Value = 19
$ws.Range("AG14").Value = 12

